$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.274.14'
$ws.Range('E2').Value = '  +3.72%  '
$ws.Range('D3').Value = '1.607.90'
$ws.Range('E3').Value = '  +2.49%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '213.17'
$ws.Range('E5').Value = '  +2.62%  '
$ws.Range('E6').Value = '  -0.18%  '
$ws.Range('E7').Value = '  +2.05%  '
$ws.Range('E8').Value = '  +2.69%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.0619'
$ws.Range('E9').Value = '  +1.88%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '18.08'
$ws.Range('E10').Value = '  +1.01%  '
$ws.Range('E11').Value = '  +4.87%  '
$ws.Range('D12').Value = '1.834.45'
$ws.Range('E12').Value = '  +2.74%  '
$ws.Range('D13').Value = '1.612.26'
$ws.Range('E13').Value = '  +2.75%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.02'
$ws.Range('E14').Value = '  -0.54%  '
$ws.Range('E15').Value = '  +1.38%  '
$ws.Range('D16').Value = '26.260.91'
$ws.Range('E16').Value = '  +3.79%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '60.81'
$ws.Range('E17').Value = '  +1.93%  '
$ws.Range('D18').Value = '0.0₃0725'
$ws.Range('E18').Value = '  +2.18%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '209.35'
$ws.Range('E19').Value = '  +12.87%  '
$ws.Range('E20').Value = '  -0.25%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.25'
$ws.Range('E21').Value = '  +2.75%  '
$ws.Range('E22').Value = '  +0.34%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '6.03'
$ws.Range('E23').Value = '  +2.41%  '
$ws.Range('E24').Value = '  +8.93%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '142.43'
$ws.Range('E25').Value = '  +0.97%  '
$ws.Range('E26').Value = '  -0.25%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.124'
$ws.Range('E27').Value = '  -4.00%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '15.29'
$ws.Range('E28').Value = '  +2.91%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.47'
$ws.Range('E29').Value = '  +0.16%  '
$ws.Range('E30').Value = '  +1.37%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0472'
$ws.Range('E31').Value = '  +1.94%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.14'
$ws.Range('E32').Value = '  +2.93%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.01'
$ws.Range('E33').Value = '  +0.48%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.48'
$ws.Range('E34').Value = '  +1.38%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.36'
$ws.Range('E35').Value = '  +2.24%  '
$ws.Range('D36').Value = '1.109.27'
$ws.Range('E36').Value = '  +1.68%  '
$ws.Range('E37').Value = '  +7.29%  '
$ws.Range('E38').Value = '  +0.07%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.34'
$ws.Range('E39').Value = '  +0.80%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.785'
$ws.Range('E40').Value = '  +1.16%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.498'
$ws.Range('E41').Value = '  +0.32%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.777'
$ws.Range('E42').Value = '  +1.54%  '
$ws.Range('D43').Value = '1.747.31'
$ws.Range('E43').Value = '  +2.76%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '93.09'
$ws.Range('E44').Value = '  +0.50%  '
$ws.Range('E45').Value = '  +0.98%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = '0.0₆0106'
$ws.Range('E46').Value = '  -2.79%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.55'
$ws.Range('E47').Value = '  +9.33%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '53.57'
$ws.Range('E48').Value = '  +1.44%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0506'
$ws.Range('E49').Value = '  +0.25%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.410'
$ws.Range('E51').Value = '  -0.09%  '
